$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "248.90"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "15"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "22.57"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "15"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.283"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "15"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05685"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "15"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.412"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "15"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.337"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "15"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8068"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "15"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8957"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "15"
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.01106"
$ws.Range("E10").Value = "9OneONE"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "15"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1426"
$ws.Range("E11").Value = "10WazirXWRX"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "15"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07444"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "15"
$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03093"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "15"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03004"
$ws.Range("E14").Value = "13BitrueCoinBTR"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "15"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09421"
$ws.Range("E15").Value = "14BitMartTokenBMX"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "15"
$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.871"
$ws.Range("E16").Value = "15MCDexMCB"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "15"
$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.001589"
$ws.Range("E17").Value = "16BitForexTokenBF"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "15"
$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.04799"
$ws.Range("E18").Value = "17CoinExTokenCET"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "15"
$ws.Range("B19").Value = "UpBots"
$ws.Range("C19").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.01828"
$ws.Range("E19").Value = "18UpBotsUBXTBestin24h"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "15"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.006416"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "15"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.004986"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "15"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0009979"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "15"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.0001501"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "15"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "15"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.172"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "15"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "15"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1353"
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "15"
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "15"
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "15"
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "15"
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "15"
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "15"
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "15"
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "15"
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "15"
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "15"
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "15"
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "15"
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "15"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03995"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = "15"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006836"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = "15"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1070"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = "15"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002794"
$ws.Range("E43").Value = "42CEJICEJIWorstin24h"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = "15"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007719"
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = "15"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005592"
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "15"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000751"
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = "15"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4996"
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = "15"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.2033"
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = "15"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002102"
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = "15"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.01011"
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = "15"
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = "15"
